$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.568.37"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.915.62"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("E4").Value = "  +0.10%  "
$c = $ws.Range("D5")
$c.Value = "'315.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("E6").Value = "  -0.03%  "
$c = $ws.Range("D7")
$c.Value = "'0.5157"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.58%  "
$c = $ws.Range("D8")
$c.Value = "'0.3989"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "
$c = $ws.Range("D9")
$c.Value = "'0.09805"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("E10").Value = "  +2.86%  "
$c = $ws.Range("D11")
$c.Value = "'42.29"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.01%  "
$c = $ws.Range("D12")
$c.Value = "'6.513"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.56%  "
$c = $ws.Range("D13")
$c.Value = "'21.19"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "1.915.60"
$ws.Range("E14").Value = "  +2.92%  "
$c = $ws.Range("D15")
$c.Value = "'7.464"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D17")
$c.Value = "'94.77"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D18")
$c.Value = "'0.00001137"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$c = $ws.Range("D19")
$c.Value = "'0.06660"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "
$c = $ws.Range("D20")
$c.Value = "'18.23"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.98%  "
$c = $ws.Range("D21")
$c.Value = "'1.0000"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$c = $ws.Range("D22")
$c.Value = "'6.301"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").Value = "28.630.16"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D26")
$c.Value = "'2.689"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +8.33%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.137.00"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.Value = "'21.30"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D29")
$c.Value = "'157.47"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D30")
$c.Value = "'129.77"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D31")
$c.Value = "'1.120"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +6.86%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D32")
$c.Value = "'0.1075"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D33")
$c.Value = "'5.731"
$c.Style = "Normal"
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D34")
$c.Value = "'3.634"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D35")
$c.Value = "'9.864"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +7.39%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D36")
$c.Value = "'0.06774"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D37")
$c.Value = "'0.02439"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D38")
$c.Value = "'1.275"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.01%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D39")
$c.Value = "'0.2233"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D40")
$c.Value = "'11.81"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D41")
$c.Value = "'0.6486"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D42")
$c.Value = "'5.086"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D43")
$c.Value = "'1.190"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Range("D44")
$c.Value = "'0.9998"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D45")
$c.Value = "'13.54"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D46")
$c.Value = "'0.6101"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D47")
$c.Value = "'3.779"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D48")
$c.Value = "'1.287"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D49")
$c.Value = "'2.071"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.33%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D50")
$c.Value = "'124.89"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$c = $ws.Range("D51")
$c.Value = "'1.206"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.52%  "
